# Generate Report for Handoff
# Update the "Latest Handoff Date"/"Latest Handoff Datetime" values for the
# 43cfe7d3-9445-4ff9-a125-ac988621c1dc row across the Overview, zh-cn and
# de-de worksheets to reflect a new handoff.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D5").Value = "2016-18-09 09:18:02"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E5").Value = "2016-03-09 09:17:54"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E5").Value = "2016-03-09 09:18:02"
